# "Generate Report for Handback" - localization-status.xlsx
#
# A handback has now completed for a.md (zh-cn and de-de). This script
# updates the generated localization-status report to reflect that:
#   - Status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it appears
#     (Overview sheet + each language sheet's Status column).
#   - Each language sheet's "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns get populated for a.md, which is
#     now the only row that has been handed back.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3061528868bf82d47e4311f1f8f2ab9c52a62dbf/e2e/"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns for both rows (a.md, b.md)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Range("E1:F3").Columns.AutoFit()

# ---------------------------------------------------------------------
# Per-language sheets: zh-cn and de-de
# ---------------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; HandbackTime = "2016-08-18 12:36:41" },
    @{ Sheet = "de-de"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; HandbackTime = "2016-08-18 12:36:50" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status column (C) for both rows
    $ws.Range("C2").Value = $statusNew
    $ws.Range("C3").Value = $statusNew

    # Latest Target File (I) / Latest Handback File (J) / Latest Handback DateTime (K)
    # for the a.md row (row 2) and the b.md row (row 3) - only a.md was
    # actually handed back, but the report regenerates the target-file
    # hyperlink + handoff artifacts for every localized row.
    $ws.Range("J2").Value = $lang.Xlf
    $ws.Range("K2").Value = $lang.HandbackTime

    $ws.Range("J3").Value = $lang.Xlf
    $ws.Range("K3").Value = $lang.HandbackTime

    # Add the new "Latest Target File" hyperlinks. The existing A2/A3
    # (a.md / b.md) hyperlinks are left untouched.
    $ws.Hyperlinks.Add($ws.Range("I2"), $baseUrl + "a.md", "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $baseUrl + "a.md", "", "", "a.md")

    $ws.Columns.Item(3).AutoFit()
    $ws.Columns.Item(10).AutoFit()
}
